$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: append a new (separate) run containing a single space
# right after the "Het converteren ... bereiken" run, inside the
# same paragraph.
# ------------------------------------------------------------------
$doelText = "Het converteren van een RGBImage naar een IntensityImage. Dit willen we zo goed en snel mogelijk bereiken"
$r1 = $d.Content
$found1 = $r1.Find.Execute($doelText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the 'Het converteren...' paragraph text"
}
$insertPos = $r1.End

# Splitting the paragraph and re-joining it keeps the appended text in
# its own <w:r>, instead of it being silently coalesced into the
# preceding run.
$breakRange = $d.Range($insertPos, $insertPos)
$breakRange.InsertParagraphAfter()
$newRunRange = $d.Range($insertPos + 1, $insertPos + 1)
$newRunRange.InsertAfter(" ")
$joinRange = $d.Range($insertPos, $insertPos + 1)
$joinRange.Delete()

# ------------------------------------------------------------------
# Change 2: drop the _GoBack bookmark that currently sits after
# "Dit is de methode die al was geschreven."
# ------------------------------------------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# ------------------------------------------------------------------
# Change 3: re-create the _GoBack bookmark inside the word "stukken"
# (splitting it into "st" | "ukken"), inside the Werkwijze paragraph.
# ------------------------------------------------------------------
$r3 = $d.Content
$found3 = $r3.Find.Execute("stukken", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find 'stukken'"
}
$splitPos = $r3.Start + 2
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)

# Re-stamp the text of the run that follows the new bookmark so the
# engine drops the (harmless but unneeded) xml:space="preserve" left
# over from splitting the original, space-leading run.
$tailText = "ukken code geschreven voor het converteren van een RGBImage naar een IntensityImage. Deze hebben we door een timer in te schakelen kunnen meten op snelheid."
$r4 = $d.Content
$found4 = $r4.Find.Execute($tailText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $txt = $r4.Text
    $r4.Delete()
    $r4.InsertAfter($txt)
}

Write-Output "edit complete"
